$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'272.99"
$ws.Range("E2").Value = "'4.41%"
$ws.Range("D3").Value = "'26.91"
$ws.Range("E3").Value = "'-0.89%"
$ws.Range("D4").Value = "'4.714"
$ws.Range("E4").Value = "'0.21%"
$ws.Range("D5").Value = "'0.06142"
$ws.Range("E5").Value = "'-1.06%"
$ws.Range("D6").Value = "'6.745"
$ws.Range("E6").Value = "'0.27%"
$ws.Range("D7").Value = "'0.8607"
$ws.Range("E7").Value = "'1.10%"
$ws.Range("D8").Value = "'0.9140"
$ws.Range("E8").Value = "'0.70%"
$ws.Range("D9").Value = "'0.1445"
$ws.Range("E9").Value = "'2.95%"
$ws.Range("D10").Value = "'0.05300"
$ws.Range("E10").Value = "'12.09%"
$ws.Range("D11").Value = "'0.07156"
$ws.Range("E11").Value = "'0.86%"
$ws.Range("D12").Value = "'0.03180"
$ws.Range("E12").Value = "'0.17%"
$ws.Range("D13").Value = "'0.09046"
$ws.Range("E13").Value = "'-0.17%"
$ws.Range("D14").Value = "'0.001531"
$ws.Range("E14").Value = "'-0.55%"
$ws.Range("D15").Value = "'0.0006084"
$ws.Range("D16").Value = "'0.005936"
$ws.Range("E16").Value = "'-1.43%"
$ws.Range("D17").Value = "'3.471"
$ws.Range("E17").Value = "'0.09%"
$ws.Range("D18").Value = "'3.184"
$ws.Range("E18").Value = "'0.41%"
$ws.Range("E19").Value = "'3.98%"
$ws.Range("E20").Value = "'-0.69%"
$ws.Range("D21").Value = "'0.1307"
$ws.Range("E21").Value = "'1.20%"
$ws.Range("D22").Value = "'3.836"
$ws.Range("E22").Value = "'-6.73%"
$ws.Range("D23").Value = "'0.04245"
$ws.Range("E23").Value = "'0.59%"
$ws.Range("D24").Value = "'0.001178"
$ws.Range("E24").Value = "'-3.35%"
$ws.Range("D25").Value = "'0.004195"
$ws.Range("E25").Value = "'1.93%"
$ws.Range("D26").Value = "'0.0001199"
$ws.Range("E26").Value = "'-0.13%"
$ws.Range("D27").Value = "'0.0001673"
$ws.Range("E27").Value = "'3.51%"
$ws.Range("D40").Value = "'0.03973"
$ws.Range("E40").Value = "'1.89%"
$ws.Range("D41").Value = "'0.006198"
$ws.Range("E41").Value = "'50.00%"
$ws.Range("D42").Value = "'0.1129"
$ws.Range("E42").Value = "'1.57%"
$ws.Range("D43").Value = "'0.002168"
$ws.Range("E43").Value = "'-0.73%"
$ws.Range("D44").Value = "'0.01274"
$ws.Range("E44").Value = "'-5.21%"
$ws.Range("D45").Value = "'0.00005130"
$ws.Range("E45").Value = "'-0.88%"
$ws.Range("D46").Value = "'0.00000000749"
$ws.Range("E46").Value = "'-0.30%"
$ws.Range("D47").Value = "'0.8978"
$ws.Range("E47").Value = "'465.00%"
$ws.Range("D48").Value = "'0.02986"
$ws.Range("E48").Value = "'-14.80%"
$ws.Range("D49").Value = "'0.00002096"
$ws.Range("E49").Value = "'-0.30%"
$ws.Range("D50").Value = "'0.0001996"
$ws.Range("E50").Value = "'-0.30%"
